$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(33, 8).Value = 762.9
$ws.Cells.Item(33, 10).Value = 799.25
$ws.Cells.Item(33, 12).Value = 799.25
$ws.Cells.Item(33, 14).Value = -1257.25

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(55, 8).Value = 59.444443
$ws.Cells.Item(55, 10).Value = 67.8
$ws.Cells.Item(55, 12).Value = 67.8
$ws.Cells.Item(55, 14).Value = -495.8

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(70, 8).Value = 1310.2
$ws.Cells.Item(70, 9).Value = 934
$ws.Cells.Item(70, 11).Value = 2802
$ws.Cells.Item(70, 13).Value = -2532

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(73, 8).Value = 1310.2
$ws.Cells.Item(73, 9).Value = 934
$ws.Cells.Item(73, 11).Value = 2802
$ws.Cells.Item(73, 13).Value = -1866

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(99, 8).Value = 111114936
$ws.Cells.Item(99, 9).Value = 987
$ws.Cells.Item(99, 11).Value = 2961
$ws.Cells.Item(99, 13).Value = -1463

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(100, 8).Value = 2340
$ws.Cells.Item(100, 9).Value = 2325.25
$ws.Cells.Item(100, 10).Value = 2399
$ws.Cells.Item(100, 11).Value = 2325.25
$ws.Cells.Item(100, 12).Value = 2399
$ws.Cells.Item(100, 13).Value = -1784.25
$ws.Cells.Item(100, 14).Value = -3481

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(101, 8).Value = 2548.2222
$ws.Cells.Item(101, 9).Value = 2548.2222
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 11).Value = 7644.6666
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(101, 13).Value = -6022.6666
$ws.Cells.Item(101, 14).ClearContents()

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(113, 8).Value = 2321.6667
$ws.Cells.Item(113, 9).Value = 1641
$ws.Cells.Item(113, 11).Value = 1641
$ws.Cells.Item(113, 13).Value = 1613

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(132, 8).Value = 343415.3
$ws.Cells.Item(132, 9).Value = 381955.38
$ws.Cells.Item(132, 10).Value = 15824.75
$ws.Cells.Item(132, 11).Value = 1145866.14
$ws.Cells.Item(132, 12).Value = 47474.25
$ws.Cells.Item(132, 13).Value = -1143336.14
$ws.Cells.Item(132, 14).Value = -52534.25

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(135, 8).Value = 4290.9116
$ws.Cells.Item(135, 9).Value = 1933.1333
$ws.Cells.Item(135, 11).Value = 17398.1997
$ws.Cells.Item(135, 13).Value = -14863.1997

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(137, 8).Value = 9876.719999999999
$ws.Cells.Item(137, 9).Value = 6669
$ws.Cells.Item(137, 11).Value = 20007
$ws.Cells.Item(137, 13).Value = -17457

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(138, 8).Value = 4119.8433
$ws.Cells.Item(138, 9).Value = 1181.4286
$ws.Cells.Item(138, 10).Value = 6176.7334
$ws.Cells.Item(138, 11).Value = 3544.2858
$ws.Cells.Item(138, 12).Value = 18530.2002
$ws.Cells.Item(138, 13).Value = 1595.7142
$ws.Cells.Item(138, 14).Value = -28810.2002

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(141, 8).Value = 857.4583
$ws.Cells.Item(141, 9).Value = 737.7317
$ws.Cells.Item(141, 11).Value = 2213.1951
$ws.Cells.Item(141, 13).Value = 2966.8049

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 14556751
$ws.Cells.Item(32, 9).Value = 15156134
$ws.Cells.Item(32, 11).Value = 15156134
$ws.Cells.Item(32, 13).Value = -15155847

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 3157.8262
$ws.Cells.Item(61, 9).Value = 2737.513
$ws.Cells.Item(61, 11).Value = 2737.513
$ws.Cells.Item(61, 13).Value = -2525.513

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(74, 8).Value = 5630.381
$ws.Cells.Item(74, 9).Value = 6305
$ws.Cells.Item(74, 10).Value = 4534.125
$ws.Cells.Item(74, 11).Value = 6305
$ws.Cells.Item(74, 12).Value = 4534.125
$ws.Cells.Item(74, 13).Value = -5431
$ws.Cells.Item(74, 14).Value = -6282.125

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(77, 8).Value = 5630.381
$ws.Cells.Item(77, 9).Value = 6305
$ws.Cells.Item(77, 10).Value = 4534.125
$ws.Cells.Item(77, 11).Value = 31525
$ws.Cells.Item(77, 12).Value = 22670.625
$ws.Cells.Item(77, 13).Value = -27157
$ws.Cells.Item(77, 14).Value = -31406.625

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(97, 8).Value = 964.7037
$ws.Cells.Item(97, 9).Value = 753.7778
$ws.Cells.Item(97, 10).Value = 1386.5555
$ws.Cells.Item(97, 11).Value = 753.7778
$ws.Cells.Item(97, 12).Value = 1386.5555
$ws.Cells.Item(97, 13).Value = -257.7778
$ws.Cells.Item(97, 14).Value = -2378.5555

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(134, 8).Value = 30214.5
$ws.Cells.Item(134, 10).Value = 30214.5
$ws.Cells.Item(134, 12).Value = 30214.5
$ws.Cells.Item(134, 14).Value = -40354.5

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(136, 8).Value = 3157.8262
$ws.Cells.Item(136, 9).Value = 2737.513
$ws.Cells.Item(136, 11).Value = 8212.539000000001
$ws.Cells.Item(136, 13).Value = -5662.539000000001

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(80, 8).Value = 14504795
$ws.Cells.Item(80, 9).Value = 1927.3334
$ws.Cells.Item(80, 10).Value = 23828066
$ws.Cells.Item(80, 11).Value = 1927.3334
$ws.Cells.Item(80, 12).Value = 23828066
$ws.Cells.Item(80, 13).Value = -929.3334
$ws.Cells.Item(80, 14).Value = -23830062

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(83, 8).Value = 14504795
$ws.Cells.Item(83, 9).Value = 1927.3334
$ws.Cells.Item(83, 10).Value = 23828066
$ws.Cells.Item(83, 11).Value = 9636.666999999999
$ws.Cells.Item(83, 12).Value = 119140330
$ws.Cells.Item(83, 13).Value = -4644.666999999999
$ws.Cells.Item(83, 14).Value = -119150314

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(86, 8).Value = 28998.666
$ws.Cells.Item(86, 9).Value = 100000
$ws.Cells.Item(86, 11).Value = 100000
$ws.Cells.Item(86, 13).Value = -98877

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(89, 8).Value = 28998.666
$ws.Cells.Item(89, 9).Value = 100000
$ws.Cells.Item(89, 11).Value = 500000
$ws.Cells.Item(89, 13).Value = -494384

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 1396653.8
$ws.Cells.Item(134, 9).Value = 1569516.8
$ws.Cells.Item(134, 11).Value = 4708550.4
$ws.Cells.Item(134, 13).Value = -4706015.4

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(62, 8).Value = 11333.667
$ws.Cells.Item(62, 9).Value = 27500
$ws.Cells.Item(62, 10).Value = 3250.5
$ws.Cells.Item(62, 11).Value = 27500
$ws.Cells.Item(62, 12).Value = 3250.5
$ws.Cells.Item(62, 13).Value = -26876
$ws.Cells.Item(62, 14).Value = -4498.5

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(65, 8).Value = 11333.667
$ws.Cells.Item(65, 9).Value = 27500
$ws.Cells.Item(65, 10).Value = 3250.5
$ws.Cells.Item(65, 11).Value = 137500
$ws.Cells.Item(65, 12).Value = 16252.5
$ws.Cells.Item(65, 13).Value = -134380
$ws.Cells.Item(65, 14).Value = -22492.5

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(94, 8).Value = 52636428
$ws.Cells.Item(94, 9).Value = 90912600
$ws.Cells.Item(94, 10).Value = 6690.75
$ws.Cells.Item(94, 11).Value = 90912600
$ws.Cells.Item(94, 12).Value = 6690.75
$ws.Cells.Item(94, 13).Value = -90912149
$ws.Cells.Item(94, 14).Value = -7592.75

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(122, 8).Value = 1919.9736
$ws.Cells.Item(122, 9).Value = 1471.0834
$ws.Cells.Item(122, 10).Value = 10000
$ws.Cells.Item(122, 11).Value = 4413.2502
$ws.Cells.Item(122, 12).Value = 30000
$ws.Cells.Item(122, 13).Value = -1963.2502
$ws.Cells.Item(122, 14).Value = -34900

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(132, 8).Value = 7284
$ws.Cells.Item(132, 9).Value = 6372.9
$ws.Cells.Item(132, 10).Value = 9308.666999999999
$ws.Cells.Item(132, 11).Value = 19118.7
$ws.Cells.Item(132, 12).Value = 27926.001
$ws.Cells.Item(132, 13).Value = -16588.7
$ws.Cells.Item(132, 14).Value = -32986.001

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(134, 8).Value = 62511630
$ws.Cells.Item(134, 9).Value = 142866820
$ws.Cells.Item(134, 11).Value = 428600460
$ws.Cells.Item(134, 13).Value = -428597925

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(5, 8).Value = 313448.7
$ws.Cells.Item(5, 9).Value = 714.2857
$ws.Cells.Item(5, 10).Value = 556686.5600000001
$ws.Cells.Item(5, 11).Value = 2142.8571
$ws.Cells.Item(5, 12).Value = 1670059.68
$ws.Cells.Item(5, 13).Value = -2030.8571
$ws.Cells.Item(5, 14).Value = -1670283.68

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(37, 8).Value = 97918.13
$ws.Cells.Item(37, 10).Value = 97918.13
$ws.Cells.Item(37, 12).Value = 293754.39
$ws.Cells.Item(37, 14).Value = -293978.39

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(135, 8).Value = 313448.7
$ws.Cells.Item(135, 9).Value = 714.2857
$ws.Cells.Item(135, 10).Value = 556686.5600000001
$ws.Cells.Item(135, 11).Value = 6428.571300000001
$ws.Cells.Item(135, 12).Value = 5010179.040000001
$ws.Cells.Item(135, 13).Value = -3893.571300000001
$ws.Cells.Item(135, 14).Value = -5015249.040000001

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(107, 8).Value = 941.0909
$ws.Cells.Item(107, 10).Value = 999.2
$ws.Cells.Item(107, 12).Value = 999.2
$ws.Cells.Item(107, 14).Value = -4839.2

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(122, 8).Value = 35338.29
$ws.Cells.Item(122, 9).Value = 54706.26
$ws.Cells.Item(122, 10).Value = 4672.3335
$ws.Cells.Item(122, 11).Value = 164118.78
$ws.Cells.Item(122, 12).Value = 14017.0005
$ws.Cells.Item(122, 13).Value = -161668.78
$ws.Cells.Item(122, 14).Value = -18917.0005

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(46, 8).Value = 20001630
$ws.Cells.Item(46, 10).Value = 31252004
$ws.Cells.Item(46, 12).Value = 31252004
$ws.Cells.Item(46, 14).Value = -31252380

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(93, 8).Value = 2926.5483
$ws.Cells.Item(93, 9).Value = 2476.842
$ws.Cells.Item(93, 11).Value = 2476.842
$ws.Cells.Item(93, 13).Value = -1228.842

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(132, 8).Value = 2741.1365
$ws.Cells.Item(132, 9).Value = 2837.5
$ws.Cells.Item(132, 11).Value = 8512.5
$ws.Cells.Item(132, 13).Value = -5982.5

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(62, 8).Value = 13312.5
$ws.Cells.Item(62, 9).Value = 11150
$ws.Cells.Item(62, 11).Value = 11150
$ws.Cells.Item(62, 13).Value = -10526

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(65, 8).Value = 13312.5
$ws.Cells.Item(65, 9).Value = 11150
$ws.Cells.Item(65, 11).Value = 55750
$ws.Cells.Item(65, 13).Value = -52630

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value = 1921.2142
$ws.Cells.Item(81, 9).Value = 1755.5555
$ws.Cells.Item(81, 11).Value = 3511.111
$ws.Cells.Item(81, 13).Value = -2450.111

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(84, 8).Value = 1921.2142
$ws.Cells.Item(84, 9).Value = 1755.5555
$ws.Cells.Item(84, 11).Value = 17555.555
$ws.Cells.Item(84, 13).Value = -12251.555

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(96, 8).Value = 2495
$ws.Cells.Item(96, 10).Value = 2495
$ws.Cells.Item(96, 12).Value = 2495
$ws.Cells.Item(96, 14).Value = -5241

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132, 8).Value = 6410.0967
$ws.Cells.Item(132, 9).Value = 5464.077
$ws.Cells.Item(132, 11).Value = 16392.231
$ws.Cells.Item(132, 13).Value = -13862.231
